$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date and Publisher values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-parameter"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": update the embedded ibm.com -> linuxforhealth.org URLs ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-parameter"
$elements.Range("J6").Value = "ParameterDefinition {http://linuxforhealth.org/fhir/cdm/StructureDefinition/parameter-definition-with-default}
"
